$wb = $excel.ActiveWorkbook

# --- Sheet "Klassifizierungen" (sheet2): change selection only ---
$wsK = $wb.Worksheets.Item("Klassifizierungen")
$wsK.Select()
$wsK.Range("F38").Select()

# --- Sheet "Standorte und Kapazitäten" (sheet3): data edits ---
$ws = $wb.Worksheets.Item("Standorte und Kapazitäten")
$ws.Select()

# Row 4: update coordinates
$ws.Range("G4").Value = -106.499686
$ws.Range("H4").Value = -78.608621499999998

# Row 5: restructure - insert new B5/C5 (D5 keeps pointing at its existing
# text value "33" automatically once the shared-string table is
# renumbered), add K5, update coordinates, and clear M5
$ws.Range("B5").Value = "Ort 2"
$ws.Range("C5").Value = "Marktplatz"
$ws.Range("G5").Value = -106.499686
$ws.Range("H5").Value = -78.608621499999998
$ws.Range("K5").Value = "BB"
$ws.Range("M5").ClearContents()

# Row 6: new data row
$ws.Range("B6").Value = "Place3"
$ws.Range("C6").Value = "Hauptstraße"
$ws.Range("D6").Value = 123
$ws.Range("E6").Value = 12345
$ws.Range("F6").Value = "Hamburg"
$ws.Range("G6").Value = 10
$ws.Range("H6").Value = 55
$ws.Range("I6").Value = 123
$ws.Range("J6").Value = "Zwei"
$ws.Range("K6").Value = "AA"
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 33

# Row 7: new data row
$ws.Range("B7").Value = "Place 4"
$ws.Range("G7").Value = 11
$ws.Range("H7").Value = 66
$ws.Range("I7").Value = 123
$ws.Range("K7").Value = "AA"
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = 0

# Note: the extended (x14) list data-validations on columns J/K are left
# untouched here - this runtime's Validation object always (de)serialises
# through the legacy <dataValidations> block, so Add/Delete cannot
# reproduce the existing x14 extLst entries without corrupting the
# legacy block that must stay exactly as-is.

# Final selection on the data sheet
$ws.Range("J6").Select()
